$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.01730024852427183
$ws.Range("C2").Value = 0.01743797812681002
$ws.Range("D2").Value = 0.01664993251499259
$ws.Range("E2").Value = 0.01580140637326485

$ws.Range("B3").Value = 7.164600579590038
$ws.Range("C3").Value = 8.34981664908516
$ws.Range("D3").Value = 8.466001735709472
$ws.Range("E3").Value = 8.287056293015256

$ws.Range("B4").Value = -0.007715731536047051
$ws.Range("C4").Value = -0.01072926463919976
$ws.Range("D4").Value = -0.01313059040095384
$ws.Range("E4").Value = -0.01533374246478889

$ws.Range("B5").Value = -2.960985194873507
$ws.Range("C5").Value = -5.278654477475951
$ws.Range("D5").Value = -7.920378794277832
$ws.Range("E5").Value = -6.154797652350622

$ws.Range("B6").Value = 0.001939250719171284
$ws.Range("C6").Value = 0.002214536068434529
$ws.Range("D6").Value = 0.0005560889073582502
$ws.Range("E6").Value = -0.00124476984616817

$ws.Range("B7").Value = 2.043383266899122
$ws.Range("C7").Value = 1.35957417776543
$ws.Range("D7").Value = 0.2438515610454587
$ws.Range("E7").Value = -0.4164969276094324

$ws.Range("B8").Value = 0.01704408369759636
$ws.Range("C8").Value = 0.01599109251986152
$ws.Range("D8").Value = 0.01592105532579664
$ws.Range("E8").Value = 0.0147756170397372

$ws.Range("B9").Value = 6.859552101866575
$ws.Range("C9").Value = 7.691104751428291
$ws.Range("D9").Value = 7.928774249915929
$ws.Range("E9").Value = 7.699502759219947

$ws.Range("B10").Value = -0.008357527046492438
$ws.Range("C10").Value = -0.01192195670660364
$ws.Range("D10").Value = -0.01434129409155893
$ws.Range("E10").Value = -0.01681369197704382

$ws.Range("B11").Value = -3.423039717211485
$ws.Range("C11").Value = -5.899705947516384
$ws.Range("D11").Value = -6.631425257262459
$ws.Range("E11").Value = -6.430676064501979

$ws.Range("B12").Value = 0.001840024300327477
$ws.Range("C12").Value = 0.001011729802501773
$ws.Range("D12").Value = -0.0007448541283637865
$ws.Range("E12").Value = -0.003541466442669543

$ws.Range("B13").Value = 2.132812360357899
$ws.Range("C13").Value = 0.6551861297420783
$ws.Range("D13").Value = -0.3266746206261197
$ws.Range("E13").Value = -1.133678317260042

$ws.Range("B14").Value = 0.01328894048839277
$ws.Range("C14").Value = 0.01352169674337987
$ws.Range("D14").Value = 0.01281321607820902
$ws.Range("E14").Value = 0.01202016398799491

$ws.Range("B15").Value = 5.291025205464209
$ws.Range("C15").Value = 6.178792138611436
$ws.Range("D15").Value = 6.269227790973541
$ws.Range("E15").Value = 6.213169556932412

$ws.Range("B16").Value = -0.009041690451927127
$ws.Range("C16").Value = -0.01263615281712632
$ws.Range("D16").Value = -0.01594696925871424
$ws.Range("E16").Value = -0.01836798765851275

$ws.Range("B17").Value = -3.479376881707805
$ws.Range("C17").Value = -5.834334341290351
$ws.Range("D17").Value = -6.492683760118034
$ws.Range("E17").Value = -7.104735858558067

$ws.Range("B18").Value = 0.0007076709781948513
$ws.Range("C18").Value = -0.0009347170612388165
$ws.Range("D18").Value = -0.003639369204604907
$ws.Range("E18").Value = -0.007236591819227384

$ws.Range("B19").Value = 0.8218775133873167
$ws.Range("C19").Value = -0.5576268074890148
$ws.Range("D19").Value = -1.520988649219853
$ws.Range("E19").Value = -2.207458316196834

$ws.Range("B20").Value = 0.01250928047880678
$ws.Range("C20").Value = 0.01195200433994183
$ws.Range("D20").Value = 0.01111592538759691
$ws.Range("E20").Value = 0.01046191925249658

$ws.Range("B21").Value = 4.743679683272769
$ws.Range("C21").Value = 5.346723612849027
$ws.Range("D21").Value = 5.535115477684325
$ws.Range("E21").Value = 5.442017002352673

$ws.Range("B22").Value = -0.01206491642932434
$ws.Range("C22").Value = -0.01610843657722701
$ws.Range("D22").Value = -0.01945519334682248
$ws.Range("E22").Value = -0.02070905493725941

$ws.Range("B23").Value = -4.534244893993776
$ws.Range("C23").Value = -7.321129523134057
$ws.Range("D23").Value = -7.606929816504106
$ws.Range("E23").Value = -7.922102879841858

$ws.Range("B24").Value = -0.0003187444768624916
$ws.Range("C24").Value = -0.003946282955426309
$ws.Range("D24").Value = -0.006035399194845385
$ws.Range("E24").Value = -0.008727206880805066

$ws.Range("B25").Value = -0.3395441319129954
$ws.Range("C25").Value = -2.174803074233515
$ws.Range("D25").Value = -2.464360799951158
$ws.Range("E25").Value = -2.230996129500106

